$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "AUT S.R.I." column (F) entirely.
$ws.Range("F1:F1").EntireColumn.Delete() | Out-Null

# Remove the trailing "PAGO", "ENTREGA", "CAMBIO" columns (originally N:P,
# now shifted to M:O after the first deletion).
$ws.Range("M1:O1").EntireColumn.Delete() | Out-Null

# Restore the cursor/selection position recorded in the saved file.
$ws.Range("H20").Select() | Out-Null
